# Refresh the GSC "Breadcrumbs" export for heatlabs.net: append the newest
# day (2025-11-02) to the "Chart" sheet's trend table.
#
# The "Chart" sheet stores its Date column (A) as plain text, e.g. "2025-11-01",
# not as a real Excel date serial. If we just assign a date-shaped string to
# Range.Value, Excel's smart-text parser will silently convert it to a date
# serial number (and mint a date-formatted style for it), which would not
# match the existing column. To keep the new cell a literal text value (same
# as every other cell in column A) we temporarily force a text number format
# before the assignment, then clear the format again so the cell ends up back
# on the sheet's default (unstyled) look, exactly like its neighbors.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Next free row right after the last populated row of the Date column (A).
$xlUp = -4162
$lastRow = $ws.Cells(1048576, 1).End($xlUp).Row
$newRow = $lastRow + 1

$dateCell = $ws.Range("A" + $newRow)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-11-02"
$dateCell.ClearFormats()

$ws.Range("B" + $newRow).Value = 0
$ws.Range("C" + $newRow).Value = 115
